$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Aula 32 - "Fragmentando o cabecalho e rodape" - new observation entry
# (adds a second row for lesson 32, same session "7. Thymeleaf para as Views")
$ws.Range("B22").Value = 32
$ws.Range("C22").Value = "7. Thymeleaf para as Views"
$ws.Range("D22").Value = "32 . Fragmentando o cabeçalho e rodapé"
$ws.Range("D22").WrapText = $true
$ws.Range("E22").Value = "1:28`n7. Thymeleaf para as Views`n32. Fragmentando o cabeçalho e rodapé`npara usar os componentes do thymeleaf na IDE, é interessante habilitar o plugin para que a IDE reconheça e faça o auto preenchimento. Para habilitar utilizando o a ide STS basta clicar direito com direito em cima do projeto, ir em Thymeleaf e clicar em Add Thymeleaf Nature"
$ws.Range("E22").WrapText = $true
$ws.Rows.Item(22).RowHeight = 105

# Aula 33 - "Fragmentando o sidebar e a pagina home"
$ws.Range("B23").Value = 33
$ws.Range("C23").Value = "7. Thymeleaf para as Views"
$ws.Range("D23").Value = "33. Fragmentando o sidebar e a página home"
$ws.Range("E23").Value = "2:36 - aplicando fragmento através da biblioteca de layout o ""layout:fragment"" invés do th:fragment do thymeleaf. a biblioteca layout:fragment fornece recursos mais avançados."
$ws.Range("E23").WrapText = $true
$ws.Rows.Item(23).RowHeight = 45

# Match the author's final selection/scroll position
$ws.Range("D23").Select() | Out-Null
